{"js": "// Office.js (Word JavaScript API) edit script.\n// The document contains a single table of \"two-digit \u00f7 one-digit\" division\n// prompts (one prompt per non-empty cell, blank spacer rows in between).\n// The commit replaces each existing prompt with a new one, in document\n// (row-major, left-to-right / top-to-bottom) order. Because a few source\n// and/or target strings repeat (e.g. \"85\u00f74=\" appears twice, and \"33\u00f72=\"\n// is both a source and a later target), we must apply the replacements\n// positionally \u2014 NOT via a naive global text -> text dictionary lookup.\n\n// Ordered (old -> new) pairs exactly as they occur, top-to-bottom /\n// left-to-right, walking only the non-blank cells.\nconst replacements = [\n  [\"25\u00f78=\", \"29\u00f77=\"],\n  [\"65\u00f79=\", \"90\u00f76=\"],\n  [\"34\u00f76=\", \"60\u00f75=\"],\n  [\"91\u00f73=\", \"84\u00f77=\"],\n  [\"62\u00f77=\", \"84\u00f75=\"],\n  [\"46\u00f73=\", \"70\u00f76=\"],\n  [\"68\u00f75=\", \"33\u00f72=\"],\n  [\"81\u00f75=\", \"36\u00f77=\"],\n  [\"67\u00f77=\", \"59\u00f76=\"],\n  [\"59\u00f78=\", \"83\u00f72=\"],\n  [\"53\u00f78=\", \"62\u00f78=\"],\n  [\"95\u00f76=\", \"22\u00f74=\"],\n  [\"35\u00f79=\", \"84\u00f78=\"],\n  [\"85\u00f74=\", \"76\u00f79=\"],\n  [\"59\u00f77=\", \"84\u00f74=\"],\n  [\"85\u00f74=\", \"46\u00f77=\"],\n  [\"27\u00f73=\", \"76\u00f78=\"],\n  [\"61\u00f78=\", \"12\u00f74=\"],\n  [\"74\u00f72=\", \"61\u00f73=\"],\n  [\"33\u00f77=\", \"64\u00f77=\"],\n  [\"33\u00f72=\", \"36\u00f76=\"],\n  [\"98\u00f77=\", \"10\u00f75=\"],\n  [\"41\u00f72=\", \"98\u00f72=\"],\n  [\"62\u00f75=\", \"70\u00f73=\"],\n  [\"51\u00f75=\", \"40\u00f74=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst values = table.values;\nlet idx = 0;\nconst newValues = values.map((row) =>\n  row.map((cell) => {\n    if (cell === \"\" || cell === null || cell === undefined) {\n      return cell;\n    }\n    if (idx < replacements.length) {\n      const [oldVal, newVal] = replacements[idx];\n      idx++;\n      if (cell !== oldVal) {\n        throw new Error(\n          `Unexpected cell text at position ${idx - 1}: expected \"${oldVal}\" but found \"${cell}\"`\n        );\n      }\n      return newVal;\n    }\n    return cell;\n  })\n);\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# The document contains a single table of \"two-digit \u00f7 one-digit\" division\n# prompts (one prompt per non-empty cell, with blank spacer rows in\n# between). The commit replaces each existing prompt with a new one, in\n# document (row-major, left-to-right / top-to-bottom) order. Because a\n# few source and/or target strings repeat (e.g. \"85\u00f74=\" appears twice,\n# and \"33\u00f72=\" is both a source and a later target), we apply the\n# replacements positionally as we walk the table cells, instead of doing\n# a naive global text -> text Find/Replace.\n\n$replacements = @(\n    @(\"25\u00f78=\", \"29\u00f77=\"),\n    @(\"65\u00f79=\", \"90\u00f76=\"),\n    @(\"34\u00f76=\", \"60\u00f75=\"),\n    @(\"91\u00f73=\", \"84\u00f77=\"),\n    @(\"62\u00f77=\", \"84\u00f75=\"),\n    @(\"46\u00f73=\", \"70\u00f76=\"),\n    @(\"68\u00f75=\", \"33\u00f72=\"),\n    @(\"81\u00f75=\", \"36\u00f77=\"),\n    @(\"67\u00f77=\", \"59\u00f76=\"),\n    @(\"59\u00f78=\", \"83\u00f72=\"),\n    @(\"53\u00f78=\", \"62\u00f78=\"),\n    @(\"95\u00f76=\", \"22\u00f74=\"),\n    @(\"35\u00f79=\", \"84\u00f78=\"),\n    @(\"85\u00f74=\", \"76\u00f79=\"),\n    @(\"59\u00f77=\", \"84\u00f74=\"),\n    @(\"85\u00f74=\", \"46\u00f77=\"),\n    @(\"27\u00f73=\", \"76\u00f78=\"),\n    @(\"61\u00f78=\", \"12\u00f74=\"),\n    @(\"74\u00f72=\", \"61\u00f73=\"),\n    @(\"33\u00f77=\", \"64\u00f77=\"),\n    @(\"33\u00f72=\", \"36\u00f76=\"),\n    @(\"98\u00f77=\", \"10\u00f75=\"),\n    @(\"41\u00f72=\", \"98\u00f72=\"),\n    @(\"62\u00f75=\", \"70\u00f73=\"),\n    @(\"51\u00f75=\", \"40\u00f74=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$idx = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        # A cell's Range.Text includes a trailing end-of-cell marker\n        # character; strip it (and any trailing CR) so we compare/assign\n        # just the visible text.\n        $raw = $cell.Range.Text\n        $txt = $raw.TrimEnd([char]7, [char]13)\n        if ($txt -ne \"\") {\n            $pair = $replacements[$idx]\n            $oldVal = $pair[0]\n            $newVal = $pair[1]\n            if ($txt -ne $oldVal) {\n                throw (\"Unexpected cell text at row {0}, col {1}: expected '{2}' but found '{3}'\" -f $r, $c, $oldVal, $txt)\n            }\n            $cell.Range.Text = $newVal\n            $idx++\n        }\n    }\n}\n"}
